$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column (Price) values are plain text cells ("inlineStr") in the source data.
# Excel auto-converts numeric-looking strings (e.g. "547.20") to numbers when
# assigned directly, which would lose the original text formatting/precision.
# Force text via NumberFormat "@", then restore the (unstyled) look by copying the
# neighboring Link cell style, so no stray cell-level style attribute is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.647.79"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.24"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.20"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.97"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.305.60"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.88"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.718.18"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.599.62"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.324.86"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.67"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.40"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.47"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.25"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.07"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.31"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.81"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.78"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.96"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.08"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "297.32"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.31"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0950"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0501"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.43"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.62"

# Restore default (unstyled) look for the D cells touched above
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("D50").Style = $ws.Range("C50").Style

# E column (Volume 1h) values are already non-numeric-looking text ("  -2.27%  ")
# so a direct .Value assignment keeps them as text without any style changes.

$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("E3").Value = "  -4.14%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("E6").Value = "  -2.64%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("E15").Value = "  -4.16%  "
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("E18").Value = "  -3.12%  "
$ws.Range("E19").Value = "  -4.40%  "
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("E22").Value = "  -4.48%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("E25").Value = "  -6.94%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -6.08%  "
$ws.Range("E28").Value = "  -5.98%  "
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  -5.22%  "
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("E33").Value = "  -5.52%  "
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -5.56%  "
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("E42").Value = "  -7.90%  "
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("E48").Value = "  -6.88%  "
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("E51").Value = "  -0.44%  "
